$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 479, shifting rows 479:583 down to 480:584
$ws.Rows.Item(479).Insert()

# Populate the new row 479 with the required values
$ws.Cells.Item(479, 1).Value = 5
$ws.Cells.Item(479, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(479, 3).Value = "Maule"
$ws.Cells.Item(479, 4).Value = 45258
$ws.Cells.Item(479, 5).Value = 7
$ws.Cells.Item(479, 6).Value = 100112006
$ws.Cells.Item(479, 7).Value = "Repollo"
$ws.Cells.Item(479, 8).Value = "Crespo record"
$ws.Cells.Item(479, 9).Value = "Primera"
$ws.Cells.Item(479, 10).Value = 4000
$ws.Cells.Item(479, 11).Value = 900
$ws.Cells.Item(479, 12).Value = 900
$ws.Cells.Item(479, 13).Value = 900
$ws.Cells.Item(479, 14).Value = "`$/unidad"
$ws.Cells.Item(479, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(479, 16).Value = 900
$ws.Cells.Item(479, 17).Value = 1
$ws.Cells.Item(479, 18).Value = "Hortaliza"
